# Apply the "updated ExperStrat-WSX, Platform filters test cases" edit:
#  - The StatQuery column (C) for the Participants/Samples/Files tab rows now
#    includes an additional OPTIONAL MATCH on (diag:diagnosis) and groups by it.
#  - Row heights for the affected rows shrink to fit the (now shorter) query text.
#  - The active selection moves to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatsQuery = "MATCH (s:study)<--(p:participant)`n" +
                 "OPTIONAL MATCH (p)<--(samp:sample)`n" +
                 "MATCH (samp)<--(f:file)`n" +
                 "WHERE f.experimental_strategy_and_data_subtypes in ['WXS']`n" +
                 "OPTIONAL MATCH (p)<--(diag:diagnosis)`n" +
                 "WITH DISTINCT samp,diag,s,p,f`n" +
                 "RETURN`n" +
                 "    count(distinct s) AS Studies,`n" +
                 "    count(distinct p) AS Participants,`n" +
                 "    count(distinct samp) AS Samples,`n" +
                 "    count(distinct f) AS ``Files``"

$ws.Range("C2").Value = $newStatsQuery
$ws.Range("C3").Value = $newStatsQuery
$ws.Range("C4").Value = $newStatsQuery

$ws.Rows.Item(2).RowHeight = 186
$ws.Rows.Item(3).RowHeight = 186
$ws.Rows.Item(4).RowHeight = 186

$ws.Range("B5").Select()
